$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Existing rows 2-5: "Results" column all become SKIP, and the
# "Description" text shifts up one row (since two new watch/unwatch rows
# were inserted ahead of the "SKIP" string in the shared pool) -- net
# effect is just these explicit values:
$ws.Range("E2").Value = "SKIP"
$ws.Range("E3").Value = "SKIP"
$ws.Range("E4").Value = "SKIP"
$ws.Range("E5").Value = "SKIP"

# New rows for TestCase_E5 / TestCase_E6, formatted the same as row 5
# (thin border all around; column C additionally wraps text)
$ws.Range("A5:E5").Copy()
$ws.Range("A6:E7").PasteSpecial(-4122)

$ws.Range("A6").Value = "TestCase_E5"
$ws.Range("B6").Value = "TBD-01"
$ws.Range("C6").Value = "Verify that user is able to watch an Patent from ALL content search results page"
$ws.Range("D6").Value = "Y"
$ws.Range("E6").Value = "PASS"

$ws.Range("A7").Value = "TestCase_E6"
$ws.Range("B7").Value = "TBD-02"
$ws.Range("C7").Value = "Verify that user is able to watch an Post from ALL content search results page"
$ws.Range("D7").Value = "Y"
$ws.Range("E7").Value = "PASS"

# Match the author's final cursor position (cell A7 selected)
$ws.Range("A7").Select() | Out-Null
